$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# id=13 name="Straight Connector 12"
$sh = $s.Shapes.Item(2)
$sh.Top = 242.82149606299214

# id=8 name="Graphic 7"
$sh = $s.Shapes.Item(3)
$sh.Top = 195.2247244094488

# id=47 name="Graphic 46"
$sh = $s.Shapes.Item(4)
$sh.Top = 225.42945861816406

# id=4 name="Rectangle 3"
$sh = $s.Shapes.Item(6)
$sh.Top = 130.09023622047243
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(1,1).ParagraphFormat.Alignment = 3
$tr.Paragraphs(2,1).ParagraphFormat.Alignment = 3

# id=19 name="Graphic 18"
$sh = $s.Shapes.Item(7)
$sh.Left = 458.8390808105469
$sh.Top = 346.4729309082031
$sh.Width = 28.728662490844727
$sh.Height = 28.728662490844727

# id=23 name="Graphic 22"
$sh = $s.Shapes.Item(8)
$sh.Left = 415.4864807128906
$sh.Top = 342.8255905511811

# id=53 name="Graphic 52"
$sh = $s.Shapes.Item(9)
$sh.Top = 215.76062992125983

# id=55 name="Rectangle 54"
$sh = $s.Shapes.Item(10)
$sh.Top = 130.09023622047243

# id=56 name="Graphic 55"
$sh = $s.Shapes.Item(11)
$sh.Top = 195.2247244094488

# id=58 name="Straight Connector 57"
$sh = $s.Shapes.Item(12)
$sh.Top = 279.173095703125

# id=60 name="Rectangle 59"
$sh = $s.Shapes.Item(13)
$sh.Top = 376.522705078125

# id=61 name="Rectangle 60"
$sh = $s.Shapes.Item(14)
$sh.Top = 219.47291338582679

# id=63 name="Straight Connector 62"
$sh = $s.Shapes.Item(15)
$sh.Top = 242.82149606299214

# id=64 name="Rectangle 63"
$sh = $s.Shapes.Item(16)
$sh.Top = 206.46992125984252

# id=65 name="Straight Connector 64"
$sh = $s.Shapes.Item(17)
$sh.Top = 269.88244094488186

# id=66 name="Graphic 65"
$sh = $s.Shapes.Item(18)
$sh.Top = 218.1251968503937

# id=68 name="Graphic 67"
$sh = $s.Shapes.Item(19)
$sh.Top = 340.18341064453125

# id=69 name="Rectangle 68"
$sh = $s.Shapes.Item(20)
$sh.Top = 328.4332580566406
$sh.Height = 60.58590551181102
$sh.TextFrame.TextRange.Paragraphs(2,1).Font.Size = 10

# id=70 name="Straight Connector 69"
$sh = $s.Shapes.Item(21)
$sh.Top = 269.88244094488186
